$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

# Pre-format the LAST UPDATE column as text so the date-like string
# "04-Nov-2025" is stored literally instead of being auto-converted
# into a date serial number by Excel's smart-entry parsing.
$dateRange = $ws.Range("I3:I31")
$dateRange.NumberFormat = "@"

for ($row = 3; $row -le 31; $row++) {
    $hCell = $ws.Cells.Item($row, 8)   # Column H - PERIOD TO EXPIRE
    $hCell.Value2 = $hCell.Value2 - 1

    $iCell = $ws.Cells.Item($row, 9)   # Column I - LAST UPDATE
    $iCell.Value = "04-Nov-2025"
}
